# Update countries & provincias Spain
# Applies the data refresh captured in the commit:
#  - Timestamp in A1 moves from 08:14 to 09:31
#  - Estados Unidos (row 4), India (row 6), Singapur (row 46) and
#    Hungria (row 100) get refreshed totals
#  - Armenia now ranks ahead of Honduras/Afganistan (rows 52-54) with
#    fresh numbers, Honduras/Afganistan keep their previous figures but
#    shift down a row
#  - Lituania now ranks ahead of Guinea-Bisau (rows 126-127) the same way
#  - Georgia now ranks ahead of Uganda (rows 144-145) the same way

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 09:31"

# --- Straight numeric refreshes ---------------------------------------
# Estados Unidos
$ws.Range("B4").Value = 4101000
$ws.Range("C4").Value = 125
$ws.Range("D4").Value = 1942815
$ws.Range("E4").Value = 2012000
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 146185

# India
$ws.Range("B6").Value = 1241416
$ws.Range("C6").Value = 1732
$ws.Range("D6").Value = 784432
$ws.Range("E6").Value = 427080
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 29904

# Singapur
$ws.Range("B46").Value = 49098
$ws.Range("C46").Value = 354
$ws.Range("E46").Value = 4276

# Hungria
$ws.Range("B100").Value = 4380
$ws.Range("C100").Value = 14
$ws.Range("D100").Value = 3300
$ws.Range("E100").Value = 484

# --- Armenia jumps ahead of Honduras / Afganistan (rows 52-54) --------
$ws.Range("A52").Value = "Armenia"
$ws.Range("B52").Value = 36162
$ws.Range("C52").Value = 469
$ws.Range("D52").Value = 25244
$ws.Range("E52").Value = 10230
$ws.Range("G52").Value = 10
$ws.Range("H52").Value = 688

$ws.Range("A53").Value = "Honduras"
$ws.Range("B53").Value = 36102
$ws.Range("C53").Value = 757
$ws.Range("D53").Value = 4315
$ws.Range("E53").Value = 30781
$ws.Range("G53").Value = 18
$ws.Range("H53").Value = 1006

$ws.Range("A54").Value = "Afganistan"
$ws.Range("B54").Value = 35915
$ws.Range("C54").Value = 188
$ws.Range("D54").Value = 24538
$ws.Range("E54").Value = 10166
$ws.Range("G54").Value = 21
$ws.Range("H54").Value = 1211

# --- Lituania jumps ahead of Guinea-Bisau (rows 126-127) ---------------
$ws.Range("A126").Value = "Lituania"
$ws.Range("B126").Value = 1960
$ws.Range("C126").Value = 9
$ws.Range("D126").Value = 1611
$ws.Range("E126").Value = 269
$ws.Range("H126").Value = 80

$ws.Range("A127").Value = "Guinea-Bisau"
$ws.Range("B127").Value = 1954
$ws.Range("D127").Value = 803
$ws.Range("E127").Value = 1125
$ws.Range("H127").Value = 26

# --- Georgia jumps ahead of Uganda (rows 144-145) -----------------------
$ws.Range("A144").Value = "Georgia"
$ws.Range("B144").Value = 1085
$ws.Range("C144").Value = 12
$ws.Range("D144").Value = 911
$ws.Range("E144").Value = 158
$ws.Range("H144").Value = 16

$ws.Range("A145").Value = "Uganda"
$ws.Range("B145").Value = 1075
$ws.Range("D145").Value = 958
$ws.Range("E145").Value = 117
$ws.Range("H145").Value = 0
